$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - flow_base_flow
$ws.Range("B4").Value = 89
$ws.Range("D4").Value = 0.7

# Row 6 - temperature_rearing
$ws.Range("B6").Value = 54
$ws.Range("D6").Value = 0.43

# Row 7 - riparian
$ws.Range("B7").Value = 120
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0.94
$ws.Range("E7").Value = 0
